# Apply attendance updates to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: mark Invalid (G) and Absent (H)
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Rows 4-8: mark Absent (H)
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1

# Row 9: mark Total Attendance Count (D) and Real (E)
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

# Rows 10-18: mark Absent (H)
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
